$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 688
$ws.Range("J4").Value = 1000
$ws.Range("L4").Value = 1000
$ws.Range("N4").Value = -1228
$ws.Range("H5").Value = 385.18182
$ws.Range("I5").Value = 418.625
$ws.Range("J5").Value = 296
$ws.Range("K5").Value = 418.625
$ws.Range("L5").Value = 296
$ws.Range("M5").Value = -303.625
$ws.Range("N5").Value = -526
$ws.Range("H28").Value = 481.375
$ws.Range("I28").Value = 72
$ws.Range("K28").Value = 72
$ws.Range("M28").Value = 413
$ws.Range("H33").Value = 424.3684
$ws.Range("I33").Value = 392.5
$ws.Range("K33").Value = 392.5
$ws.Range("M33").Value = -163.5
$ws.Range("H74").Value = 2579
$ws.Range("I74").Value = 723.75
$ws.Range("K74").Value = 723.75
$ws.Range("M74").Value = 212.25
$ws.Range("H77").Value = 2579
$ws.Range("I77").Value = 723.75
$ws.Range("K77").Value = 3618.75
$ws.Range("M77").Value = 1061.25
$ws.Range("H92").Value = 439.3846
$ws.Range("I92").Value = 439.3846
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 439.3846
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("M92").Value = 808.6154

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 203.42857
$ws.Range("J5").Value = 60
$ws.Range("L5").Value = 60
$ws.Range("N5").Value = -284
$ws.Range("H102").Value = 1760.5555
$ws.Range("I102").Value = 1320.7142
$ws.Range("J102").Value = 3300
$ws.Range("K102").Value = 1320.7142
$ws.Range("L102").Value = 3300
$ws.Range("M102").Value = 301.2858000000001
$ws.Range("N102").Value = -6544
$ws.Range("H110").Value = 1770.5714
$ws.Range("I110").Value = 952.5
$ws.Range("J110").Value = 2861.3333
$ws.Range("K110").Value = 952.5
$ws.Range("L110").Value = 2861.3333
$ws.Range("M110").Value = 1092.5
$ws.Range("N110").Value = -6951.3333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 203.42857
$ws.Range("J4").Value = 60
$ws.Range("L4").Value = 60
$ws.Range("N4").Value = -290
$ws.Range("H95").Value = 18208
$ws.Range("J95").Value = 18208
$ws.Range("L95").Value = 18208
$ws.Range("N95").Value = -23700
$ws.Range("H96").Value = 2700
$ws.Range("I96").Value = 2700
$ws.Range("K96").Value = 2700
$ws.Range("M96").Value = 46
$ws.Range("H99").Value = 1949.75
$ws.Range("I99").Value = 1949.75
$ws.Range("K99").Value = 1949.75
$ws.Range("M99").Value = -451.75
$ws.Range("H105").Value = 2675.1428
$ws.Range("I105").Value = 2675.1428
$ws.Range("K105").Value = 2675.1428
$ws.Range("M105").Value = -928.1428000000001
$ws.Range("H107").Value = 3625.3333
$ws.Range("I107").Value = 3625.3333
$ws.Range("K107").Value = 3625.3333
$ws.Range("M107").Value = -1705.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 164.33333
$ws.Range("I7").Value = 176.28572
$ws.Range("J7").Value = 122.5
$ws.Range("K7").Value = 176.28572
$ws.Range("L7").Value = 122.5
$ws.Range("M7").Value = -63.28572
$ws.Range("N7").Value = -348.5
$ws.Range("H16").Value = 1089.3
$ws.Range("I16").Value = 866
$ws.Range("J16").Value = 1982.5
$ws.Range("K16").Value = 866
$ws.Range("L16").Value = 1982.5
$ws.Range("M16").Value = -579
$ws.Range("N16").Value = -2556.5
$ws.Range("H105").Value = 458.16666
$ws.Range("I105").Value = 458.16666
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 458.16666
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()
$ws.Range("H107").Value = 815
$ws.Range("I107").Value = 700
$ws.Range("J107").Value = 1083.3334
$ws.Range("K107").Value = 700
$ws.Range("L107").Value = 1083.3334
$ws.Range("M107").Value = 1220
$ws.Range("N107").Value = -4923.3334
$ws.Range("H113").Value = 1089.3
$ws.Range("I113").Value = 866
$ws.Range("J113").Value = 1982.5
$ws.Range("K113").Value = 866
$ws.Range("L113").Value = 1982.5
$ws.Range("M113").Value = 1304
$ws.Range("N113").Value = -6322.5
$ws.Range("H132").Value = 1623.4
$ws.Range("I132").Value = 1550.9474
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 4652.8422
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -2122.8422
$ws.Range("N132").Value = -14060
$ws.Range("M105").Value = 1288.83334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 160.74074
$ws.Range("J2").Value = 228.41667
$ws.Range("L2").Value = 1370.50002
$ws.Range("N2").Value = -1596.50002
$ws.Range("H4").Value = 1391.96
$ws.Range("I4").Value = 1304.6364
$ws.Range("K4").Value = 3913.9092
$ws.Range("M4").Value = -3801.9092
$ws.Range("H113").Value = 513.7143
$ws.Range("J113").Value = 418.6
$ws.Range("L113").Value = 1255.8
$ws.Range("N113").Value = -5595.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2188.9
$ws.Range("I80").Value = 2024
$ws.Range("J80").Value = 2573.6667
$ws.Range("K80").Value = 2024
$ws.Range("L80").Value = 2573.6667
$ws.Range("M80").Value = -1026
$ws.Range("N80").Value = -4569.6667
$ws.Range("H83").Value = 2188.9
$ws.Range("I83").Value = 2024
$ws.Range("J83").Value = 2573.6667
$ws.Range("K83").Value = 10120
$ws.Range("L83").Value = 12868.3335
$ws.Range("M83").Value = -5128
$ws.Range("N83").Value = -22852.3335
$ws.Range("H97").Value = 1070
$ws.Range("I97").Value = 1025
$ws.Range("J97").Value = 1250
$ws.Range("K97").Value = 1025
$ws.Range("L97").Value = 1250
$ws.Range("M97").Value = -529
$ws.Range("N97").Value = -2242
$ws.Range("H113").Value = 667.5
$ws.Range("I113").Value = 557.1667
$ws.Range("J113").Value = 998.5
$ws.Range("K113").Value = 557.1667
$ws.Range("L113").Value = 998.5
$ws.Range("M113").Value = 1612.8333
$ws.Range("N113").Value = -5338.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1972.6666
$ws.Range("I16").Value = 1967.2
$ws.Range("K16").Value = 1967.2
$ws.Range("M16").Value = -1797.2
$ws.Range("H22").Value = 1097.8572
$ws.Range("I22").Value = 528.3333
$ws.Range("J22").Value = 1525
$ws.Range("K22").Value = 528.3333
$ws.Range("L22").Value = 1525
$ws.Range("M22").Value = -233.3333
$ws.Range("N22").Value = -2115
$ws.Range("H27").Value = 1097.8572
$ws.Range("I27").Value = 528.3333
$ws.Range("J27").Value = 1525
$ws.Range("K27").Value = 528.3333
$ws.Range("L27").Value = 1525
$ws.Range("M27").Value = -421.3333
$ws.Range("N27").Value = -1739
$ws.Range("H61").Value = 4088.5
$ws.Range("J61").Value = 3752.5
$ws.Range("L61").Value = 3752.5
$ws.Range("N61").Value = -4156.5
$ws.Range("H82").Value = 2992.6875
$ws.Range("I82").Value = 1845.5714
$ws.Range("K82").Value = 1845.5714
$ws.Range("M82").Value = -1484.5714
$ws.Range("H85").Value = 2992.6875
$ws.Range("I85").Value = 1845.5714
$ws.Range("K85").Value = 1845.5714
$ws.Range("M85").Value = -597.5714
$ws.Range("H93").Value = 1453.7059
$ws.Range("I93").Value = 1429.0646
$ws.Range("J93").Value = 1708.3334
$ws.Range("K93").Value = 1429.0646
$ws.Range("L93").Value = 1708.3334
$ws.Range("M93").Value = -181.0645999999999
$ws.Range("N93").Value = -4204.3334
$ws.Range("H100").Value = 3400
$ws.Range("I100").Value = 3666.6667
$ws.Range("J100").Value = 3000
$ws.Range("K100").Value = 3666.6667
$ws.Range("L100").Value = 3000
$ws.Range("M100").Value = -3125.6667
$ws.Range("N100").Value = -4082
$ws.Range("H113").Value = 4088.5
$ws.Range("J113").Value = 3752.5
$ws.Range("L113").Value = 3752.5
$ws.Range("N113").Value = -8092.5
$ws.Range("H132").Value = 3682.111
$ws.Range("I132").Value = 4082.5
$ws.Range("K132").Value = 12247.5
$ws.Range("M132").Value = -9717.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1749.5
$ws.Range("I81").Value = 1749.5
$ws.Range("K81").Value = 3499
$ws.Range("M81").Value = -2438
$ws.Range("H84").Value = 1749.5
$ws.Range("I84").Value = 1749.5
$ws.Range("K84").Value = 17495
$ws.Range("M84").Value = -12191
$ws.Range("H100").Value = 884.2222
$ws.Range("I100").Value = 851.1429000000001
$ws.Range("K100").Value = 1702.2858
$ws.Range("M100").Value = -1161.2858
$ws.Range("H107").Value = 504.18182
$ws.Range("I107").Value = 341.5
$ws.Range("J107").Value = 699.4
$ws.Range("K107").Value = 1024.5
$ws.Range("L107").Value = 2098.2
$ws.Range("M107").Value = 895.5
$ws.Range("N107").Value = -5938.2
$ws.Range("H113").Value = 592.6667
$ws.Range("I113").Value = 475
$ws.Range("J113").Value = 651.5
$ws.Range("K113").Value = 1425
$ws.Range("L113").Value = 1954.5
$ws.Range("M113").Value = 745
$ws.Range("N113").Value = -6294.5
$ws.Range("H124").Value = 9421.25
$ws.Range("J124").Value = 9421.25
$ws.Range("L124").Value = 9421.25
$ws.Range("N124").Value = -19241.25
$ws.Range("H132").Value = 1831.9166
$ws.Range("I132").Value = 1648.3
$ws.Range("K132").Value = 4944.9
$ws.Range("M132").Value = -2414.9
$ws.Range("H136").Value = 1255
$ws.Range("I136").Value = 1255
$ws.Range("K136").Value = 3765
$ws.Range("M136").Value = -1215
